$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2243461.009
$ws.Range("B3").Value = 1497.819
$ws.Range("B4").Value = 1272.797
